$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) is always stored as text in this workbook, even for
# cells whose content looks like a plain decimal number (e.g. "224.28").
# Force a Text number format before writing so Excel does not silently
# reinterpret/re-normalize these strings as floating point numbers (which
# would e.g. turn "40.80" into 40.8 or drop exact decimal precision).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

$ws.Range("D2").Value = "34.415.89"
$ws.Range("E2").Value = "  +0.39%  "
$ws.Range("D3").Value = "1.798.58"
$ws.Range("E3").Value = "  +0.36%  "
$ws.Range("D4").Value = "1.01"
$ws.Range("E4").Value = "  +0.39%  "
$ws.Range("D5").Value = "224.28"
$ws.Range("E5").Value = "  -0.23%  "
$ws.Range("D6").Value = "0.601"
$ws.Range("E6").Value = "  +0.98%  "
$ws.Range("E7").Value = "  +0.35%  "
$ws.Range("D8").Value = "40.80"
$ws.Range("E8").Value = "  +12.81%  "
$ws.Range("D9").Value = "0.290"
$ws.Range("E9").Value = "  +0.11%  "
$ws.Range("E10").Value = "  -1.02%  "
$ws.Range("D11").Value = "0.1000"
$ws.Range("E11").Value = "  +3.94%  "
$ws.Range("D12").Value = "2.058.11"
$ws.Range("E12").Value = "  +0.36%  "
$ws.Range("D13").Value = "1.794.51"
$ws.Range("E13").Value = "  -0.32%  "
$ws.Range("D14").Value = "10.78"
$ws.Range("E14").Value = "  -2.97%  "
$ws.Range("D15").Value = "34.407.19"
$ws.Range("E15").Value = "  +0.39%  "
$ws.Range("E16").Value = "  -0.72%  "
$ws.Range("D17").Value = "4.38"
$ws.Range("E17").Value = "  +0.53%  "
$ws.Range("D18").Value = "67.18"
$ws.Range("E18").Value = "  -1.96%  "
$ws.Range("D19").Value = "239.24"
$ws.Range("E19").Value = "  -0.23%  "
$ws.Range("E20").Value = "  -0.53%  "
$ws.Range("D21").Value = "11.07"
$ws.Range("E21").Value = "  -1.31%  "
$ws.Range("E22").Value = "  +0.32%  "
$ws.Range("D23").Value = "4.10"
$ws.Range("E23").Value = "  +1.15%  "
$ws.Range("E24").Value = "  -0.81%  "
$ws.Range("D25").Value = "171.57"
$ws.Range("E25").Value = "  +0.68%  "
$ws.Range("D26").Value = "7.62"
$ws.Range("E26").Value = "  -4.41%  "
$ws.Range("D27").Value = "17.26"
$ws.Range("E27").Value = "  +0.73%  "
$ws.Range("E28").Value = "  +0.46%  "
$ws.Range("D29").Value = "1.01"
$ws.Range("E29").Value = "  +0.41%  "
$ws.Range("E30").Value = "  +0.16%  "
$ws.Range("E31").Value = "  +0.25%  "
$ws.Range("E32").Value = "  -0.70%  "
$ws.Range("D33").Value = "0.0509"
$ws.Range("E33").Value = "  -0.43%  "
$ws.Range("D34").Value = "1.76"
$ws.Range("E34").Value = "  +0.65%  "
$ws.Range("D35").Value = "1.315.63"
$ws.Range("E35").Value = "  -3.01%  "
$ws.Range("D36").Value = "0.643"
$ws.Range("E36").Value = "  +0.12%  "
$ws.Range("E37").Value = "  +0.87%  "
$ws.Range("D38").Value = "85.48"
$ws.Range("E38").Value = "  +6.58%  "
$ws.Range("E39").Value = "  +1.66%  "
$ws.Range("E40").Value = "  -0.46%  "
$ws.Range("D41").Value = "14.66"
$ws.Range("E41").Value = "  +11.72%  "
$ws.Range("E42").Value = "  +6.63%  "
$ws.Range("E43").Value = "  +0.94%  "
$ws.Range("D44").Value = "2.79"
$ws.Range("E44").Value = "  +0.27%  "
$ws.Range("D45").Value = "0.935"
$ws.Range("E45").Value = "  +0.61%  "
$ws.Range("D46").Value = "0.0518"
$ws.Range("D47").Value = "1.959.77"
$ws.Range("E47").Value = "  +0.41%  "
$ws.Range("E48").Value = "  +1.53%  "
$ws.Range("E49").Value = "  +0.30%  "
$ws.Range("D50").Value = "100.40"
$ws.Range("E50").Value = "  -1.20%  "
$ws.Range("E51").Value = "  +1.48%  "
